# Updates symbol list (prices, and a re-ranked block of coins rows 10-18)
# to match the "Updated symbol list" GitHub Actions commit.
# Values are text (stored as inlineStr in the source), so we prefix with a
# leading apostrophe to force Excel to keep them as text instead of
# auto-converting numeric-looking strings into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'244.57"
$ws.Range("D3").Value = "'21.96"
$ws.Range("D4").Value = "'5.447"
$ws.Range("D5").Value = "'0.05768"
$ws.Range("D6").Value = "'3.423"
$ws.Range("D7").Value = "'6.318"
$ws.Range("D9").Value = "'1.021"
$ws.Range("B10").Value = "'WazirX"
$ws.Range("C10").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1424"
$ws.Range("E10").Value = "'9WazirXWRX"
$ws.Range("B11").Value = "'MandalaExchangeToken"
$ws.Range("C11").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07300"
$ws.Range("E11").Value = "'10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03136"
$ws.Range("E12").Value = "'11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "'BitrueCoin"
$ws.Range("C13").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03114"
$ws.Range("E13").Value = "'12BitrueCoinBTR"
$ws.Range("B14").Value = "'MCDex"
$ws.Range("C14").Value = "'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D14").Value = "'4.139"
$ws.Range("E14").Value = "'13MCDexMCB"
$ws.Range("B15").Value = "'BitMartToken"
$ws.Range("C15").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09377"
$ws.Range("E15").Value = "'14BitMartTokenBMX"
$ws.Range("B16").Value = "'BitForexToken"
$ws.Range("C16").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001591"
$ws.Range("E16").Value = "'15BitForexTokenBF"
$ws.Range("B17").Value = "'CoinExToken"
$ws.Range("C17").Value = "'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04802"
$ws.Range("E17").Value = "'16CoinExTokenCET"
$ws.Range("B18").Value = "'One"
$ws.Range("C18").Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0005843"
$ws.Range("E18").Value = "'17OneONE"
$ws.Range("D19").Value = "'0.006296"
$ws.Range("D20").Value = "'0.004125"
$ws.Range("D21").Value = "'0.0009906"
$ws.Range("D23").Value = "'3.750"
$ws.Range("D24").Value = "'2.179"
$ws.Range("D26").Value = "'0.1328"
$ws.Range("D27").Value = "'0.0003994"
$ws.Range("D40").Value = "'0.03861"
$ws.Range("D41").Value = "'0.006640"
$ws.Range("D43").Value = "'0.002616"
$ws.Range("D44").Value = "'0.006584"
$ws.Range("D45").Value = "'0.00005596"
